# Refresh the crypto price/volume snapshot (Price = col D, Volume(1h) = col E)
# for each coin row on the active sheet, per the scheduled GitHub Actions update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.836.18'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '1.636.20'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''215.32'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').Value = '''0.505'
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').Value = '''0.0643'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').Value = '''19.89'
$ws.Range('E10').Value = '  +2.34%  '
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('D13').Value = '1.639.10'
$ws.Range('E13').Value = '  +0.62%  '
$ws.Range('D14').Value = '1.860.44'
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('E16').Value = '  +1.65%  '
$ws.Range('D17').Value = '''63.09'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').Value = '25.836.37'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = '''194.10'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').Value = '''9.94'
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('D23').Value = '''6.16'
$ws.Range('E23').Value = '  +2.74%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('E25').Value = '  -1.65%  '
$ws.Range('D26').Value = '''139.66'
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').Value = '''0.121'
$ws.Range('E27').Value = '  -4.53%  '
$ws.Range('D28').Value = '''6.84'
$ws.Range('E28').Value = '  +1.61%  '
$ws.Range('D29').Value = '''15.50'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('E30').Value = '  +0.35%  '
$ws.Range('D31').Value = '''0.0494'
$ws.Range('E31').Value = '  +1.61%  '
$ws.Range('E32').Value = '  +1.16%  '
$ws.Range('E33').Value = '  +1.53%  '
$ws.Range('E34').Value = '  +1.68%  '
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('D36').Value = '''0.903'
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('E38').Value = '  +1.02%  '
$ws.Range('D39').Value = '1.118.01'
$ws.Range('E39').Value = '  -0.94%  '
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').Value = '''99.44'
$ws.Range('E43').Value = '  +2.25%  '
$ws.Range('D44').Value = '''0.798'
$ws.Range('E44').Value = '  +0.26%  '
$ws.Range('E45').Value = '  -3.00%  '
$ws.Range('D46').Value = '''55.55'
$ws.Range('E46').Value = '  +0.65%  '
$ws.Range('D47').Value = '''2.47'
$ws.Range('E47').Value = '  +11.04%  '
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('D50').Value = '''7.62'
$ws.Range('E50').Value = '  -0.61%  '
$ws.Range('E51').Value = '  -0.10%  '
